$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3: "divesh" -> "divesh_jain"
$ws.Range("A3").Value = "divesh_jain"

# B3: numeric 9864121 -> text "1234567"
# Set the number format to Text first so Excel stores the value as a string
# rather than re-parsing it back into a number.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1234567"

# I3: "muthu street" -> "vaichur muthiah street"
$ws.Range("I3").Value = "vaichur muthiah street"
